$p = $ppt.ActivePresentation
$s = $p.Slides.Item(5)
$shape = $s.Shapes.Item(2)
$para = $shape.TextFrame.TextRange.Paragraphs(1)
$run1 = $para.Runs(1)
$run1.Text = "HTML, PHP 7, JS, CSS, TWIG, "
$run2 = $run1.InsertAfter("Materialize")
$run2.Font.Size = 24
